# Update "想去人数" (F column) values in both the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 296
$wsExhibition.Range("F3").Value = 1143
$wsExhibition.Range("F4").Value = 2593

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F2").Value = 296
$wsAllTypes.Range("F5").Value = 1143
$wsAllTypes.Range("F6").Value = 2593
